$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-06-07"

# Update the "June (through 06-06)" label to "June (through 06-07)"
$ws.Range("A7").Value = "June (through 06-07)"

# Update 2022 (column I) value for May (row 6): 114 -> 113
$ws.Range("I6").Value = 113

# Update June row (row 7) values for columns C..I
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 16
$ws.Range("G7").Value = 39
$ws.Range("H7").Value = 28
$ws.Range("I7").Value = 21

# Update Total row (row 8) values for columns C..I
$ws.Range("C8").Value = 220
$ws.Range("D8").Value = 327
$ws.Range("E8").Value = 311
$ws.Range("G8").Value = 397
$ws.Range("H8").Value = 659
$ws.Range("I8").Value = 684
